$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New job match data for rows 2-8 (Title, Company, Location, Match Score, Matched Keywords, Posted At, Apply Link)
$rows = @(
    @("Software Engineer New", "Convey", "Chicago, IL, US USA", 14.4, "RAG, Copilot, Kinesis, Docker, Kubernetes, Git, Kafka, MongoDB, NoSQL, SQL", "2026-02-23", "https://www.indeed.com/viewjob?jk=7f0c38e749fbda51"),
    @("Senior Software Engineer New", "Convey", "Chicago, IL, US USA", 14.4, "RAG, Copilot, Kinesis, Docker, Kubernetes, Git, Kafka, MongoDB, NoSQL, SQL", "2026-02-23", "https://www.indeed.com/viewjob?jk=5a7e1e1ef2ff4bf9"),
    @("Senior Security Engineer – Cloud & Data Security", "Sigma Computing", "New York, NY, US USA", 13.3, "RAG, BigQuery, Kubernetes, CI/CD, Terraform, Snowflake, Databricks, BigQuery, Python, SQL", "2026-02-23", "https://www.indeed.com/viewjob?jk=5d4ac754da67252d"),
    @("Technical Support Engineer", "Sigma Computing", "New York, NY, US USA", 12.2, "RAG, Redshift, BigQuery, Snowflake, Databricks, BigQuery, Redshift, Python, SQL, R", "2026-02-23", "https://www.indeed.com/viewjob?jk=b94386ac1471ed5e"),
    @("Technical Support Engineer", "Sigma Computing", "San Francisco, CA, US USA", 11.1, "RAG, Redshift, BigQuery, Snowflake, BigQuery, Redshift, Python, SQL, R, Scala", "2026-02-23", "https://www.indeed.com/viewjob?jk=a01c49c40c088fe3"),
    @("Senior Software Development Engineer - Advertising Technology", "Expedia Group", "Chicago, IL, US USA", 11.1, "Data Scientist, Machine Learning Engineer, Generative AI, Databricks, PySpark, Kafka, R, Java, Scala, Optimization", "2026-02-22", "https://www.indeed.com/viewjob?jk=5437a1b0e8faf95d"),
    @("ML / AI Infrastructure Engineering Intern", "Arcline", "San Francisco, CA, US USA", 11.1, "LangChain, RAG, LLaMA, Copilot, Pinecone, FastAPI, Git, Python, R, Optimization", "2026-02-23", "https://www.indeed.com/viewjob?jk=9afaefdf5d1354cd")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]

    # Force the "Posted At" date-looking string to stay stored as plain text
    # (matches source data, which stores dates as literal strings, not date serials).
    $dateCell = $ws.Cells.Item($r, 6)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[5]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 7).Value = $row[6]
    $r = $r + 1
}
